$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-21 Monday" "2025-07-22 Tuesday"

Replace-Text "597×8=" "327×2="
Replace-Text "340×7=" "396×3="
Replace-Text "902×5=" "796×4="
Replace-Text "227×4=" "199×4="
Replace-Text "463×2=" "709×5="
Replace-Text "635×3=" "241×6="
Replace-Text "126×2=" "839×3="
Replace-Text "267×6=" "992×9="
Replace-Text "256×7=" "296×2="
Replace-Text "385×3=" "605×3="
Replace-Text "436×9=" "718×8="
Replace-Text "702×9=" "102×2="
Replace-Text "125×6=" "898×3="
Replace-Text "463×5=" "357×8="
Replace-Text "246×7=" "658×6="
Replace-Text "662×6=" "160×9="
Replace-Text "302×6=" "909×4="
Replace-Text "648×3=" "855×2="
Replace-Text "186×5=" "116×8="
Replace-Text "985×8=" "728×3="
Replace-Text "225×2=" "181×4="
Replace-Text "405×8=" "816×7="
Replace-Text "767×9=" "747×9="
Replace-Text "902×3=" "105×5="
Replace-Text "655×3=" "381×2="
